# Applies the documentation wording cleanup from the commit
# "DungeonGeneratoria paranneltu ja boss huone korjattu":
#   - drops the "(kuinka monta?)" aside
#   - drops the redundant "myös" before "aloitushuoneesta"
#   - drops the two "(ei vielä implementoitu koodiin)" notes, which are no
#     longer true now that the generator code has caught up with the docs

$d = $word.ActiveDocument

# Paragraph: "Tehdään vielä joitain reittejä (kuinka monta?) satunnaisesti
# valittuihin kartan sijainteihin. Reitit alkavat myös aloitushuoneesta.
# (ei vielä implementoitu koodiin)"
# becomes: "Tehdään vielä joitain reittejä satunnaisesti valittuihin kartan
# sijainteihin. Reitit alkavat aloitushuoneesta."

$r = $d.Content
$r.Find.Execute(" (kuinka monta?)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "", 2) | Out-Null

$r = $d.Content
$r.Find.Execute("myös aloitushuoneesta. (ei vielä implementoitu koodiin)", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "aloitushuoneesta.", 2) | Out-Null

# Paragraph: "Jokaisesta kentästä löytyy avainhuone, aloitushuone ja
# pomohuone, mutta muut ruudut määritetään seuraavien todennäköisyyksien
# avulla: (ei vielä implementoitu koodiin)"
# becomes: "Jokaisesta kentästä löytyy avainhuone, aloitushuone ja
# pomohuone, mutta muut ruudut määritetään seuraavien todennäköisyyksien
# avulla."

$r = $d.Content
$r.Find.Execute("avulla: (ei vielä implementoitu koodiin)", $false, $false, `
    $false, $false, $false, $true, 1, $false, "avulla.", 2) | Out-Null
